$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON)
$ws.Range("B2").Value = 2.9921861609128104
$ws.Range("C2").Value = 0.42801434716386477
$ws.Range("D2").Value = 3.018134905795482
$ws.Range("E2").Value = 0.38689077983612485

# Row 3 (STR)
$ws.Range("B3").Value = 2.4520921249786043
$ws.Range("C3").Value = 0.91898135611724596
$ws.Range("D3").Value = 6.2121606916272807
$ws.Range("E3").Value = 1.1116130973260185
